$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C: rows 7-12 -> 0
$ws.Range("C7:C12").Value = 0

# Column C: rows 13-42 -> 0.85
$ws.Range("C13:C42").Value = 0.85

# Column B: rows 16-43 -> 0.03
$ws.Range("B16:B43").Value = 0.03

# Column E: rows 16-43 -> 0.1
$ws.Range("E16:E43").Value = 0.1

# Update selection/view: activate Sheet1 and select E16:E43 with active cell E16,
# and scroll so topLeftCell resets (matches target sheetView with no topLeftCell override)
$ws.Activate()
$ws.Range("E16:E43").Select()

Write-Host "edit complete"
